$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 4-6): D, F, G, H, I ---
$ws.Range("D4").Value = 536
$ws.Range("F4").Value = 567
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 6

$ws.Range("D5").Value = 534
$ws.Range("F5").Value = 567
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 6

$ws.Range("D6").Value = 53
$ws.Range("F6").Value = 567
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 6

# --- Table 2 (rows 12-14): D, F, G, H, I, J ---
$ws.Range("D12").Value = 536
$ws.Range("F12").Value = 567
$ws.Range("G12").Value = 533
$ws.Range("H12").Value = 75
$ws.Range("I12").Value = 567
$ws.Range("J12").Value = 43

$ws.Range("D13").Value = 534
$ws.Range("F13").Value = 567
$ws.Range("G13").Value = 533
$ws.Range("H13").Value = 75
$ws.Range("I13").Value = 567
$ws.Range("J13").Value = 43

$ws.Range("D14").Value = 53
$ws.Range("F14").Value = 567
$ws.Range("G14").Value = 533
$ws.Range("H14").Value = 75
$ws.Range("I14").Value = 567
$ws.Range("J14").Value = 43

# --- Table 3 (row 20): C, D, E / I, J, K / O, P, Q ---
$ws.Range("C20").Value = 1534
$ws.Range("D20").Value = 536
$ws.Range("E20").Value = 345

$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 534
$ws.Range("K20").Value = 645

$ws.Range("O20").Value = 745
$ws.Range("P20").Value = 53
$ws.Range("Q20").Value = 64
